$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.7287194209349384, 1.65323645889881, 0.7127328510149897, 6.48142807727062, 9.576116808119359)
    3 = @(1.505614041169197, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 4.371470058157054)
    4 = @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 5.488907176552729)
    5 = @(0.3464964993005633, 0.3375848360084654, 0.7127328510149897, 0.4998867070740569, 1.896700893398075)
    6 = @(0.7287194209349384, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 3.594575437922795)
    7 = @(0.7287194209349384, 1.65323645889881, 3.082599426703578, 6.48142807727062, 11.94598338380795)
    8 = @(3.182878228561681, 1.65323645889881, 3.082599426703578, 6.48142807727062, 14.40014219143469)
    9 = @(3.182878228561681, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 8.418600821238126)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 7).Value = $vals[4]
}
